$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.943.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "'1.889.55"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("D4").Value = "'1.020"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +1.78%  "
$ws.Range("D5").Value = "'335.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("E6").Value = "  +1.67%  "
$ws.Range("D7").Value = "'0.4659"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.56%  "
$ws.Range("D8").Value = "'0.3912"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.56%  "
$ws.Range("D9").Value = "'47.38"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.37%  "
$ws.Range("D10").Value = "'0.07997"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.61%  "
$ws.Range("D11").Value = "'1.013"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.21%  "
$ws.Range("D12").Value = "'21.63"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("D13").Value = "'1.893.14"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("D14").Value = "'5.942"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").Value = "'7.102"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.16%  "
$ws.Range("D16").Value = "'1.022"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.98%  "
$ws.Range("D17").Value = "'0.06772"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.44%  "
$ws.Range("D18").Value = "'87.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").Value = "'0.00001047"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("D20").Value = "'17.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.83%  "
$ws.Range("E21").Value = "  +1.67%  "
$ws.Range("D22").Value = "'27.964.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.93%  "
$ws.Range("D23").Value = "'5.488"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.46%  "
$ws.Range("D24").Value = "'10.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.64%  "
$ws.Range("D25").Value = "'2.346"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.65%  "
$ws.Range("D26").Value = "'2.111.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.07%  "
$ws.Range("D27").Value = "'159.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.26%  "
$ws.Range("D28").Value = "'19.97"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.39%  "
$ws.Range("D29").Value = "'2.068"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.33%  "
$ws.Range("D30").Value = "'5.401"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.90%  "
$ws.Range("D31").Value = "'121.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.83%  "
$ws.Range("D32").Value = "'0.9605"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").Value = "'0.09485"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.68%  "
$ws.Range("D34").Value = "'3.673"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.21%  "
$ws.Range("D35").Value = "'1.374"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.59%  "
$ws.Range("D36").Value = "'5.323"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("D37").Value = "'0.06091"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("D38").Value = "'0.02236"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.82%  "
$ws.Range("D39").Value = "'1.214"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.94%  "
$ws.Range("D40").Value = "'8.058"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.07%  "
$ws.Range("D41").Value = "'0.5948"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.03%  "
$ws.Range("D42").Value = "'0.1880"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.08%  "
$ws.Range("D43").Value = "'10.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("E44").Value = "  +2.19%  "
$ws.Range("D45").Value = "'0.5644"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("D46").Value = "'12.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.85%  "
$ws.Range("D47").Value = "'3.399"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("D48").Value = "'1.917"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.92%  "
$ws.Range("D49").Value = "'0.06913"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.33%  "
$ws.Range("D50").Value = "'113.64"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.38%  "
$ws.Range("D51").Value = "'1.066"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.28%  "
